{"js": "// Update the East-Asian and Complex-Script font fallbacks recorded on the\n// document's paragraph styles (vignettes/docx/bs_example.docx):\n//   - Normal / Heading: eastAsia font \"DejaVu Sans\" -> \"Tahoma\"\n//   - List / Caption / Index: add a Complex-Script (\"cs\") font of\n//     \"DejaVu Sans\" (these styles previously had no <w:rFonts> override at\n//     all, so their cs font fell back through docDefaults)\nconst styles = context.document.getStyles();\nstyles.load(\"items/nameLocal\");\nawait context.sync();\n\nconst byName = {};\nfor (const style of styles.items) {\n  byName[style.nameLocal] = style;\n}\n\n// rFonts/@w:eastAsia is exposed as Font.nameFarEast in the Word JS API.\nbyName[\"Normal\"].font.nameFarEast = \"Tahoma\";\nbyName[\"Heading\"].font.nameFarEast = \"Tahoma\";\n\n// rFonts/@w:cs is exposed as Font.nameBidirectional in the Word JS API.\nbyName[\"List\"].font.nameBidirectional = \"DejaVu Sans\";\nbyName[\"Caption\"].font.nameBidirectional = \"DejaVu Sans\";\nbyName[\"Index\"].font.nameBidirectional = \"DejaVu Sans\";\n\nawait context.sync();\n", "ps1": "# Update the East-Asian and Complex-Script font fallbacks recorded on the\n# document's paragraph styles (vignettes/docx/bs_example.docx):\n#   - Normal / Heading: eastAsia font \"DejaVu Sans\" -> \"Tahoma\"\n#   - List / Caption / Index: add a Complex-Script (\"cs\") font of\n#     \"DejaVu Sans\" (these styles previously had no <w:rFonts> override at\n#     all, so their cs font fell back through docDefaults)\n$d = $word.ActiveDocument\n\n# rFonts/@w:eastAsia is exposed as Font.NameFarEast in the Word object model.\n$normal = $d.Styles.Item(\"Normal\")\n$normal.Font.NameFarEast = \"Tahoma\"\n\n$heading = $d.Styles.Item(\"Heading\")\n$heading.Font.NameFarEast = \"Tahoma\"\n\n# rFonts/@w:cs is exposed as Font.NameBi (NameBidirectional) in the Word\n# object model.\n$list = $d.Styles.Item(\"List\")\n$list.Font.NameBi = \"DejaVu Sans\"\n\n$caption = $d.Styles.Item(\"Caption\")\n$caption.Font.NameBi = \"DejaVu Sans\"\n\n$index = $d.Styles.Item(\"Index\")\n$index.Font.NameBi = \"DejaVu Sans\"\n"}
